$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 152-155 ---
# Row 152: date 44238 -> 44448, volume 3200 -> 2400
$ws.Cells.Item(152, 4).Value = 44448
$ws.Cells.Item(152, 10).Value = 2400

# Row 153: date 44238 -> 44448, volume 1600 -> 1200
$ws.Cells.Item(153, 4).Value = 44448
$ws.Cells.Item(153, 10).Value = 1200

# Row 154: date 44399 -> 44238, volume 3460 -> 3200
$ws.Cells.Item(154, 4).Value = 44238
$ws.Cells.Item(154, 10).Value = 3200

# Row 155: date 44399 -> 44238 (volume unchanged, stays 1600)
$ws.Cells.Item(155, 4).Value = 44238

# --- Add new rows 156 and 157 (duplicates of the former 154/155 content) ---
# Row 156 (Primera)
$ws.Cells.Item(156, 1).Value = 8
$ws.Cells.Item(156, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(156, 3).Value = "Coquimbo"
$ws.Cells.Item(156, 4).Value = 44399
$ws.Cells.Item(156, 5).Value = 4
$ws.Cells.Item(156, 6).Value = 100114014
$ws.Cells.Item(156, 7).Value = "Betarraga"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 3460
$ws.Cells.Item(156, 11).Value = 450
$ws.Cells.Item(156, 12).Value = 500
$ws.Cells.Item(156, 13).Value = 475
$ws.Cells.Item(156, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(156, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(156, 16).Value = 158
$ws.Cells.Item(156, 17).Value = 3
$ws.Cells.Item(156, 18).Value = "Hortaliza"

# Row 157 (Segunda)
$ws.Cells.Item(157, 1).Value = 8
$ws.Cells.Item(157, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(157, 3).Value = "Coquimbo"
$ws.Cells.Item(157, 4).Value = 44399
$ws.Cells.Item(157, 5).Value = 4
$ws.Cells.Item(157, 6).Value = 100114014
$ws.Cells.Item(157, 7).Value = "Betarraga"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Segunda"
$ws.Cells.Item(157, 10).Value = 1600
$ws.Cells.Item(157, 11).Value = 350
$ws.Cells.Item(157, 12).Value = 400
$ws.Cells.Item(157, 13).Value = 375
$ws.Cells.Item(157, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(157, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(157, 16).Value = 125
$ws.Cells.Item(157, 17).Value = 3
$ws.Cells.Item(157, 18).Value = "Hortaliza"

# Match the date-format style (s="2") used by D2:D155 for the new D cells
$ws.Cells.Item(154, 4).NumberFormat = $ws.Cells.Item(153, 4).NumberFormat
$ws.Cells.Item(155, 4).NumberFormat = $ws.Cells.Item(153, 4).NumberFormat
$ws.Cells.Item(156, 4).NumberFormat = $ws.Cells.Item(153, 4).NumberFormat
$ws.Cells.Item(157, 4).NumberFormat = $ws.Cells.Item(153, 4).NumberFormat
